# export-student-template.xlsx
#
# The original "student" sheet was a 6-column template:
#   A: Ho ten | B: Ma | C: Ngay sinh | D: CMND/CCCD | E: SDT | F: Email
# with a hint row (2), a blank spacer row (3, merged A3:F3) and an
# example row (4).
#
# The new template only needs two columns - "Ma" (student code) first,
# then "Ho ten" (full name) - with just the header row, ready for the
# student-list upload / grade-board mapping feature.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the hint / spacer / example rows (2-4), keeping only the header row.
$ws.Range("A2:A4").EntireRow.Delete()

# Drop the trailing columns that are no longer part of the template
# (Ngay sinh, CMND/CCCD, SDT, Email).
$ws.Range("C1:F1").EntireColumn.Delete()

# Drop column A ("Ho ten"); this shifts column B ("Ma") into column A,
# carrying its header formatting and exact column width along with it.
$ws.Range("A1").EntireColumn.Delete()

# Re-create column B as "Ho ten": copy A1 (value + header style) across,
# then overwrite the text and restore the original (narrower) column
# width that used to belong to column A.
$ws.Range("A1").Copy($ws.Range("B1"))
$ws.Range("B1").Value = "Họ tên "
$ws.Columns(2).ColumnWidth = 25.14

# Final header row: A = "Ma", B = "Ho ten"
$ws.Range("A1").Value = "Mã"

$ws.Range("B8").Select() | Out-Null
